# Tech stack Blockschaltbild: Farbe der Rechtecke auf 00B0F0 angepasst
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# VBA-style RGB long value: R + G*256 + B*65536 -> results in OOXML srgbClr "00B0F0"
$blue = 0 + (176 * 256) + (240 * 65536)

$shapeNames = @("Rechteck 1", "Rechteck 2", "Rechteck 3", "Rechteck 5", "Rechteck 6")

foreach ($name in $shapeNames) {
    $shape = $s.Shapes.Item($name)
    $shape.Fill.Solid()
    $shape.Fill.ForeColor.RGB = $blue
}
